# Update gh-pages to output generated at 456a3b4
#
# Sheet "展览" (sheet1): row2 G -> "不可售"; several "想去人数" (F) counts bumped up.
# Sheet "演出" (sheet2): row2 F count bumped up.
# Sheet "本地生活" (sheet3): unchanged.
# Sheet "全部类型" (sheet4): row2 G -> "不可售"; several F counts bumped up (superset of sheet1+sheet2 rows).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F4").Value = 5289
$ws1.Range("F7").Value = 22
$ws1.Range("F8").Value = 599
$ws1.Range("F9").Value = 558
$ws1.Range("F10").Value = 1049
$ws1.Range("F12").Value = 1455
$ws1.Range("F13").Value = 4233
$ws1.Range("F14").Value = 434
$ws1.Range("F15").Value = 181
$ws1.Range("F16").Value = 161
$ws1.Range("F17").Value = 96
$ws1.Range("F18").Value = 3335
$ws1.Range("F19").Value = 164
$ws1.Range("F20").Value = 1079
$ws1.Range("F24").Value = 121
$ws1.Range("F25").Value = 37
$ws1.Range("F26").Value = 138
$ws1.Range("F27").Value = 71
$ws1.Range("F28").Value = 301
$ws1.Range("F29").Value = 29
$ws1.Range("F30").Value = 55
$ws1.Range("F33").Value = 24

# ---- Sheet "演出" ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 48

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F4").Value = 48
$ws4.Range("F5").Value = 5289
$ws4.Range("F8").Value = 22
$ws4.Range("F9").Value = 599
$ws4.Range("F10").Value = 558
$ws4.Range("F11").Value = 1049
$ws4.Range("F13").Value = 1455
$ws4.Range("F14").Value = 4233
$ws4.Range("F15").Value = 434
$ws4.Range("F16").Value = 181
$ws4.Range("F17").Value = 161
$ws4.Range("F18").Value = 96
$ws4.Range("F19").Value = 3335
$ws4.Range("F20").Value = 164
$ws4.Range("F21").Value = 1079
$ws4.Range("F25").Value = 121
$ws4.Range("F26").Value = 37
$ws4.Range("F27").Value = 138
$ws4.Range("F28").Value = 71
$ws4.Range("F29").Value = 301
$ws4.Range("F30").Value = 29
$ws4.Range("F31").Value = 55
$ws4.Range("F34").Value = 24
